$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column H (Absent) values for the "consolidated report"
$ws.Range("H10").Value = 1
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 0
